# Auto-generated PowerShell Excel COM-interop script
# chore: update Sheets via scheduled runner
# Refreshes market-board derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# for the affected Leve rows across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 182.46
$ws.Range("I15").Value = 182.46
$ws.Range("K15").Value = 547.38
$ws.Range("M15").Value = -378.38

# Row 17
$ws.Range("H17").Value = 1039.4667
$ws.Range("J17").Value = 1039.4667
$ws.Range("L17").Value = 3118.4001
$ws.Range("N17").Value = -3454.4001

# Row 40
$ws.Range("H40").Value = 910
$ws.Range("I40").Value = 910
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 910
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -735
$ws.Range("N40").ClearContents()

# Row 62
$ws.Range("H62").Value = 2573.9473
$ws.Range("I62").Value = 2617.9167
$ws.Range("J62").Value = 2498.5715
$ws.Range("K62").Value = 2617.9167
$ws.Range("L62").Value = 2498.5715
$ws.Range("M62").Value = -1993.9167
$ws.Range("N62").Value = -3746.5715

# Row 64
$ws.Range("H64").Value = 3816.3635
$ws.Range("I64").Value = 3868
$ws.Range("K64").Value = 3868
$ws.Range("M64").Value = -3620

# Row 65
$ws.Range("H65").Value = 2573.9473
$ws.Range("I65").Value = 2617.9167
$ws.Range("J65").Value = 2498.5715
$ws.Range("K65").Value = 13089.5835
$ws.Range("L65").Value = 12492.8575
$ws.Range("M65").Value = -9969.583500000001
$ws.Range("N65").Value = -18732.8575

# Row 67
$ws.Range("H67").Value = 3816.3635
$ws.Range("I67").Value = 3868
$ws.Range("K67").Value = 3868
$ws.Range("M67").Value = -3010

# Row 74
$ws.Range("H74").Value = 4011.111
$ws.Range("I74").Value = 4042.8572
$ws.Range("K74").Value = 4042.8572
$ws.Range("M74").Value = -3106.8572

# Row 76
$ws.Range("H76").Value = 3085.1853
$ws.Range("I76").Value = 3090.9092
$ws.Range("J76").Value = 3060
$ws.Range("K76").Value = 3090.9092
$ws.Range("L76").Value = 3060
$ws.Range("M76").Value = -2775.9092
$ws.Range("N76").Value = -3690

# Row 77
$ws.Range("H77").Value = 4011.111
$ws.Range("I77").Value = 4042.8572
$ws.Range("K77").Value = 20214.286
$ws.Range("M77").Value = -15534.286

# Row 79
$ws.Range("H79").Value = 3085.1853
$ws.Range("I79").Value = 3090.9092
$ws.Range("J79").Value = 3060
$ws.Range("K79").Value = 3090.9092
$ws.Range("L79").Value = 3060
$ws.Range("M79").Value = -1998.9092
$ws.Range("N79").Value = -5244

# Row 100
$ws.Range("H100").Value = 1791.6666
$ws.Range("I100").Value = 1840
$ws.Range("J100").Value = 1550
$ws.Range("K100").Value = 1840
$ws.Range("L100").Value = 1550
$ws.Range("M100").Value = -1299
$ws.Range("N100").Value = -2632

# Row 121
$ws.Range("H121").Value = 1762.262
$ws.Range("J121").Value = 1802.9269
$ws.Range("L121").Value = 5408.780699999999
$ws.Range("N121").Value = -8902.780699999999

# Row 125
$ws.Range("H125").Value = 1512.2174
$ws.Range("I125").Value = 904.5714
$ws.Range("K125").Value = 8141.1426
$ws.Range("M125").Value = -5681.1426

# Row 138
$ws.Range("H138").Value = 5044.222
$ws.Range("I138").Value = 3741
$ws.Range("J138").Value = 5873.5454
$ws.Range("K138").Value = 11223
$ws.Range("L138").Value = 17620.6362
$ws.Range("M138").Value = -6083
$ws.Range("N138").Value = -27900.6362

$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 10351.625
$ws.Range("J37").Value = 10351.625
$ws.Range("L37").Value = 10351.625
$ws.Range("N37").Value = -10897.625

# Row 63
$ws.Range("H63").Value = 2920
$ws.Range("I63").Value = 2138.4614
$ws.Range("K63").Value = 2138.4614
$ws.Range("M63").Value = -1452.4614

# Row 66
$ws.Range("H66").Value = 2920
$ws.Range("I66").Value = 2138.4614
$ws.Range("K66").Value = 10692.307
$ws.Range("M66").Value = -7260.307000000001

# Row 92
$ws.Range("H92").Value = 29999.5
$ws.Range("J92").Value = 29999.5
$ws.Range("L92").Value = 29999.5
$ws.Range("N92").Value = -34991.5

# Row 139
$ws.Range("H139").Value = 71715
$ws.Range("J139").Value = 71715
$ws.Range("L139").Value = 71715
$ws.Range("N139").Value = -81995

$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 86
$ws.Range("H86").Value = 1550.9474
$ws.Range("I86").Value = 1404.6
$ws.Range("J86").Value = 2099.75
$ws.Range("K86").Value = 1404.6
$ws.Range("L86").Value = 2099.75
$ws.Range("M86").Value = -281.5999999999999
$ws.Range("N86").Value = -4345.75

# Row 89
$ws.Range("H89").Value = 1550.9474
$ws.Range("I89").Value = 1404.6
$ws.Range("J89").Value = 2099.75
$ws.Range("K89").Value = 7023
$ws.Range("L89").Value = 10498.75
$ws.Range("M89").Value = -1407
$ws.Range("N89").Value = -21730.75

$ws = $wb.Worksheets.Item("CRP")
# Row 60
$ws.Range("H60").Value = 8567.166999999999
$ws.Range("J60").Value = 10413.25
$ws.Range("L60").Value = 10413.25
$ws.Range("N60").Value = -11435.25

# Row 68
$ws.Range("H68").Value = 18424.3
$ws.Range("J68").Value = 18886.111
$ws.Range("L68").Value = 18886.111
$ws.Range("N68").Value = -20384.111

# Row 71
$ws.Range("H71").Value = 18424.3
$ws.Range("J71").Value = 18886.111
$ws.Range("L71").Value = 56658.333
$ws.Range("N71").Value = -64146.333

$ws = $wb.Worksheets.Item("CUL")
# Row 139
$ws.Range("H139").Value = 1214.1613
$ws.Range("I139").Value = 1214.1613
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 3642.4839
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 1497.5161
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4581.8276
$ws.Range("I70").Value = 4351.7334
$ws.Range("J70").Value = 4828.357
$ws.Range("K70").Value = 4351.7334
$ws.Range("L70").Value = 4828.357
$ws.Range("M70").Value = -4081.7334
$ws.Range("N70").Value = -5368.357

# Row 73
$ws.Range("H73").Value = 4581.8276
$ws.Range("I73").Value = 4351.7334
$ws.Range("J73").Value = 4828.357
$ws.Range("K73").Value = 4351.7334
$ws.Range("L73").Value = 4828.357
$ws.Range("M73").Value = -3415.7334
$ws.Range("N73").Value = -6700.357

# Row 132
$ws.Range("H132").Value = 1888.0714
$ws.Range("I132").Value = 1578.4
$ws.Range("J132").Value = 2662.25
$ws.Range("K132").Value = 4735.200000000001
$ws.Range("L132").Value = 7986.75
$ws.Range("M132").Value = -2205.200000000001
$ws.Range("N132").Value = -13046.75

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1015.1539
$ws.Range("I46").Value = 1019.95
$ws.Range("J46").Value = 999.1667
$ws.Range("K46").Value = 1019.95
$ws.Range("L46").Value = 999.1667
$ws.Range("M46").Value = -831.95
$ws.Range("N46").Value = -1375.1667

# Row 55
$ws.Range("H55").Value = 171.28572
$ws.Range("I55").Value = 189.5
$ws.Range("J55").Value = 164
$ws.Range("K55").Value = 189.5
$ws.Range("L55").Value = 164
$ws.Range("M55").Value = -16.5
$ws.Range("N55").Value = -510

# Row 122
$ws.Range("H122").Value = 2089.4
$ws.Range("I122").Value = 2004.8485
$ws.Range("J122").Value = 2321.9167
$ws.Range("K122").Value = 6014.5455
$ws.Range("L122").Value = 6965.750100000001
$ws.Range("M122").Value = -3564.5455
$ws.Range("N122").Value = -11865.7501

$ws = $wb.Worksheets.Item("WVR")
# Row 22
$ws.Range("H22").Value = 2890.9546
$ws.Range("I22").Value = 601
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 601
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -308
$ws.Range("N22").Value = -3586

# Row 56
$ws.Range("H56").Value = 15000
$ws.Range("J56").Value = 15000
$ws.Range("L56").Value = 15000
$ws.Range("N56").Value = -16428

# Row 81
$ws.Range("H81").Value = 1479.6
$ws.Range("I81").Value = 1479.6
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2959.2
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1898.2
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 1479.6
$ws.Range("I84").Value = 1479.6
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -9492
$ws.Range("N84").ClearContents()

Write-Host "Applied scheduled market-data updates to 37 Leve rows across 8 sheets."